$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = '"4242424242424242"'
$ws.Range("A3").Value = '"4000056655665556"'
$ws.Range("A4").Value = '"5555555555554444"'
